# BOM.xlsx update: add a "Test Point" line item (row 31) and convert the
# "Qty x6" formula column into a proper shared formula spanning B2:B31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 31: Test Point BOM entry -------------------------------------
$ws.Cells.Item(31, 10).Value = "A106145CT-ND"         # J31 Part Number
$ws.Cells.Item(31, 5).Value  = "Test Point"            # E31 Device
$ws.Cells.Item(31, 7).Value  = "T1, T2, T3, T4, T5"    # G31 Parts
$ws.Cells.Item(31, 8).Value  = "Test point"            # H31 Description
$ws.Cells.Item(31, 1).Value  = 5                       # A31 Qty
$ws.Cells.Item(31, 3).Value  = 30                      # C31 Qty to Buy
$ws.Cells.Item(31, 6).Value  = 603                     # F31 Package
$ws.Cells.Item(31, 9).Value  = "Digikey"               # I31 Supplier

# --- Qty x6 column: shared formula across the whole table (B2:B31) --------
$ws.Range("B2:B31").Formula = "=A2*6"

# --- Move the saved cursor/selection to match the source workbook ---------
[void]$ws.Range("E34").Select()

# --- Page setup: clear the bogus paperSize="0" attribute ------------------
$ws.PageSetup.PaperSize = 0
